$d = $word.ActiveDocument

function Find-ParaIndex($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t.StartsWith($prefix)) { return $i }
    }
    return $null
}

$startIdx = Find-ParaIndex("Wear safety goggles")
$dullIdx = Find-ParaIndex("Dull solder joints")
$listStart = $d.Paragraphs.Item($startIdx).Range.Start
$listEnd = $d.Paragraphs.Item($dullIdx - 1).Range.End
$listRange = $d.Range($listStart, $listEnd)

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>TIPS:</w:t></w:r></w:p><w:p><w:r><w:t>-Wear safety goggles</w:t></w:r></w:p><w:p><w:r><w:t>-Avoid breathing in the solder fumes</w:t></w:r></w:p><w:p><w:r><w:t>-350°C-400°C are usual soldering temperatures</w:t></w:r></w:p><w:p><w:r><w:t>-"Pre tin" components/wires</w:t></w:r></w:p><w:p><w:r><w:t>-Heat the pad and the component while adding tin</w:t></w:r></w:p><w:p><w:r><w:t>-Use a helping hand or board clamp</w:t></w:r></w:p><w:p><w:r><w:t>-Pliers make handling small components easier</w:t></w:r></w:p><w:p><w:r><w:t>-Solder paste requires lower temperature to melt</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">-Low melt solder/hot air is the easiest way to remove </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> components with many pins</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$listRange.InsertXML($xml1)

$dullIdx2 = Find-ParaIndex("Dull solder joints")
$dullRange = $d.Paragraphs.Item($dullIdx2).Range
$dullRange.Delete()

$pageBreakPara = $d.Paragraphs.Item($dullIdx2)
$pageBreakStart = $pageBreakPara.Range.Start

$insertPoint = $d.Range($pageBreakStart, $pageBreakStart)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>-Dull solder joints aren''t great</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml2)
